$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Remove row 2 ("H 72") - all data below shifts up by one row
$ws.Rows(2).Delete()

# Apply updated missing-data mask / corrected values to the remaining rows
$ws.Range("D2").Value = -14.8
$ws.Range("C4").Value = 14.9
$ws.Range("D5").ClearContents()
$ws.Range("D7").Value = -14.4
$ws.Range("C8").ClearContents()
$ws.Range("C9").Value = 12
$ws.Range("D10").ClearContents()
$ws.Range("B11").Value = -19.9
$ws.Range("C11").ClearContents()
$ws.Range("D11").Value = -13.9
$ws.Range("B12").ClearContents()
$ws.Range("D13").Value = -14.7
$ws.Range("D14").ClearContents()
$ws.Range("D16").ClearContents()
$ws.Range("C20").Value = 13.5
$ws.Range("D21").Value = -14.7
$ws.Range("C22").ClearContents()
$ws.Range("C23").Value = 13.2
$ws.Range("D24").ClearContents()
$ws.Range("C25").ClearContents()
$ws.Range("B26").Value = -19.5
$ws.Range("B27").ClearContents()
$ws.Range("C28").Value = 12.2
$ws.Range("D28").Value = -13.9
$ws.Range("B31").Value = -19.5
$ws.Range("D31").ClearContents()
$ws.Range("B32").ClearContents()
$ws.Range("C32").ClearContents()
$ws.Range("D33").Value = -15.1
$ws.Range("B35").Value = -19.2
$ws.Range("B36").ClearContents()
$ws.Range("C36").Value = 14.3
$ws.Range("D36").ClearContents()
$ws.Range("B37").Value = -19.8
$ws.Range("B38").ClearContents()
$ws.Range("C38").ClearContents()
$ws.Range("C41").Value = 13.9
$ws.Range("C43").ClearContents()
$ws.Range("B45").Value = -19.7
$ws.Range("B46").ClearContents()
$ws.Range("C52").Value = 10.8
$ws.Range("B53").Value = -20.3
$ws.Range("C54").ClearContents()
$ws.Range("B56").ClearContents()
$ws.Range("D57").Value = -13.7
$ws.Range("D59").Value = -13.6
$ws.Range("D60").ClearContents()
$ws.Range("D62").ClearContents()
